$wb = $excel.ActiveWorkbook

# Sheet ALC, row 55 (Leve Item ID 5517)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 248.7
$ws.Range("I55").Value = 165.8
$ws.Range("K55").Value = 165.8
$ws.Range("M55").Value = 48.19999999999999

# Sheet ALC, row 74 (Leve Item ID 5507)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 8000
$ws.Range("I74").Value = 8000
$ws.Range("K74").Value = 8000
$ws.Range("M74").Value = -7064

# Sheet ALC, row 77 (Leve Item ID 5507)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 8000
$ws.Range("I77").Value = 8000
$ws.Range("K77").Value = 40000
$ws.Range("M77").Value = -35320

# Sheet ALC, row 87 (Leve Item ID 10651)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 54000
$ws.Range("J87").Value = 80000
$ws.Range("L87").Value = 80000
$ws.Range("N87").Value = -82496

# Sheet ALC, row 90 (Leve Item ID 10651)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H90").Value = 54000
$ws.Range("J90").Value = 80000
$ws.Range("L90").Value = 240000
$ws.Range("N90").Value = -252480

# Sheet ALC, row 98 (Leve Item ID 36237)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 965
$ws.Range("I98").Value = 965
$ws.Range("K98").Value = 965
$ws.Range("M98").Value = 533

# Sheet ALC, row 100 (Leve Item ID 19906)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 2342.8333
$ws.Range("I100").Value = 2342.8333
$ws.Range("K100").Value = 2342.8333
$ws.Range("M100").Value = -1801.8333

# Sheet ALC, row 103 (Leve Item ID 19909)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 2416.4546
$ws.Range("I103").Value = 451.58334
$ws.Range("J103").Value = 4774.3
$ws.Range("K103").Value = 1354.75002
$ws.Range("L103").Value = 14322.9
$ws.Range("M103").Value = -768.7500199999999
$ws.Range("N103").Value = -15494.9

# Sheet ALC, row 113 (Leve Item ID 27775)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 9669.1
$ws.Range("I113").Value = 9585.75
$ws.Range("K113").Value = 9585.75
$ws.Range("M113").Value = -6331.75

# Sheet ALC, row 122 (Leve Item ID 36237)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 965
$ws.Range("I122").Value = 965
$ws.Range("K122").Value = 2895
$ws.Range("M122").Value = -445

# Sheet ALC, row 137 (Leve Item ID 44013)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 957.125
$ws.Range("I137").Value = 872.5
$ws.Range("J137").Value = 1211
$ws.Range("K137").Value = 2617.5
$ws.Range("L137").Value = 3633
$ws.Range("M137").Value = -67.5
$ws.Range("N137").Value = -8733

# Sheet ALC, row 141 (Leve Item ID 44161)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 1285.4286
$ws.Range("I141").Value = 1166.3334
$ws.Range("K141").Value = 3499.0002
$ws.Range("M141").Value = 1680.9998

# Sheet ARM, row 2 (Leve Item ID 27713)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2477.889
$ws.Range("I2").Value = 1550.8572
$ws.Range("J2").Value = 5722.5
$ws.Range("K2").Value = 1550.8572
$ws.Range("L2").Value = 5722.5
$ws.Range("M2").Value = -1437.8572
$ws.Range("N2").Value = -5948.5

# Sheet ARM, row 32 (Leve Item ID 44147)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3198.923
$ws.Range("I32").Value = 2935.92
$ws.Range("K32").Value = 2935.92
$ws.Range("M32").Value = -2648.92

# Sheet ARM, row 63 (Leve Item ID 12528)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 4599.8
$ws.Range("J63").Value = 4500
$ws.Range("L63").Value = 4500
$ws.Range("N63").Value = -5872

# Sheet ARM, row 66 (Leve Item ID 12528)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 4599.8
$ws.Range("J66").Value = 4500
$ws.Range("L66").Value = 22500
$ws.Range("N66").Value = -29364

# Sheet ARM, row 102 (Leve Item ID 19945)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 6000
$ws.Range("I102").Value = 6000
$ws.Range("K102").Value = 6000
$ws.Range("M102").Value = -4378

# Sheet ARM, row 116 (Leve Item ID 27713)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 2477.889
$ws.Range("I116").Value = 1550.8572
$ws.Range("J116").Value = 5722.5
$ws.Range("K116").Value = 1550.8572
$ws.Range("L116").Value = 5722.5
$ws.Range("M116").Value = 743.1428000000001
$ws.Range("N116").Value = -10310.5

# Sheet ARM, row 138 (Leve Item ID 42350)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H138").Value = 60000
$ws.Range("J138").Value = 60000
$ws.Range("L138").Value = 60000
$ws.Range("N138").Value = -70280

# Sheet BSM, row 3 (Leve Item ID 27713)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2477.889
$ws.Range("I3").Value = 1550.8572
$ws.Range("J3").Value = 5722.5
$ws.Range("K3").Value = 1550.8572
$ws.Range("L3").Value = 5722.5
$ws.Range("M3").Value = -1436.8572
$ws.Range("N3").Value = -5950.5

# Sheet BSM, row 33 (Leve Item ID 1625)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H33").Value = 3659.3333
$ws.Range("I33").Value = 3659.3333
$ws.Range("K33").Value = 3659.3333
$ws.Range("M33").Value = -3323.3333

# Sheet BSM, row 86 (Leve Item ID 12526)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4553.8184
$ws.Range("I86").Value = 3939
$ws.Range("J86").Value = 5066.1665
$ws.Range("K86").Value = 3939
$ws.Range("L86").Value = 5066.1665
$ws.Range("M86").Value = -2816
$ws.Range("N86").Value = -7312.1665

# Sheet BSM, row 87 (Leve Item ID 11906)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H87").Value = 38938
$ws.Range("I87").Value = 30321
$ws.Range("K87").Value = 30321
$ws.Range("M87").Value = -29073

# Sheet BSM, row 89 (Leve Item ID 12526)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 4553.8184
$ws.Range("I89").Value = 3939
$ws.Range("J89").Value = 5066.1665
$ws.Range("K89").Value = 19695
$ws.Range("L89").Value = 25330.8325
$ws.Range("M89").Value = -14079
$ws.Range("N89").Value = -36562.8325

# Sheet BSM, row 90 (Leve Item ID 11906)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H90").Value = 38938
$ws.Range("I90").Value = 30321
$ws.Range("K90").Value = 90963
$ws.Range("M90").Value = -84723

# Sheet BSM, row 94 (Leve Item ID 19939)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 5000
$ws.Range("I94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("M94").Value = ""

# Sheet BSM, row 134 (Leve Item ID 43998)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1420.4166
$ws.Range("I134").Value = 1513.1818
$ws.Range("K134").Value = 4539.5454
$ws.Range("M134").Value = -2004.5454

# Sheet CRP, row 22 (Leve Item ID 5367)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 523.5
$ws.Range("I22").Value = 505.42856
$ws.Range("K22").Value = 505.42856
$ws.Range("M22").Value = -155.42856

# Sheet CRP, row 31 (Leve Item ID 44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2799
$ws.Range("I31").Value = 1998
$ws.Range("K31").Value = 1998
$ws.Range("M31").Value = -1703

# Sheet CRP, row 34 (Leve Item ID 44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2799
$ws.Range("I34").Value = 1998
$ws.Range("K34").Value = 1998
$ws.Range("M34").Value = -1796

# Sheet CRP, row 88 (Leve Item ID 10608)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H88").Value = 29833
$ws.Range("J88").Value = 29833
$ws.Range("L88").Value = 29833
$ws.Range("N88").Value = -30645

# Sheet CRP, row 91 (Leve Item ID 10608)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H91").Value = 29833
$ws.Range("J91").Value = 29833
$ws.Range("L91").Value = 29833
$ws.Range("N91").Value = -32641

# Sheet CUL, row 23 (Leve Item ID 4858)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 255
$ws.Range("I23").Value = 96.666664
$ws.Range("J23").Value = 492.5
$ws.Range("K23").Value = 289.999992
$ws.Range("L23").Value = 1477.5
$ws.Range("M23").Value = -54.99999200000002
$ws.Range("N23").Value = -1947.5

# Sheet CUL, row 113 (Leve Item ID 27843)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 2663.3333
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").Value = ""

# Sheet GSM, row 102 (Leve Item ID 36169)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2397.0908
$ws.Range("I102").Value = 2208.6667
$ws.Range("J102").Value = 3245
$ws.Range("K102").Value = 2208.6667
$ws.Range("L102").Value = 3245
$ws.Range("M102").Value = -586.6667000000002
$ws.Range("N102").Value = -6489

# Sheet GSM, row 122 (Leve Item ID 36182)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1820.3125
$ws.Range("I122").Value = 1625.3077
$ws.Range("K122").Value = 4875.9231
$ws.Range("M122").Value = -2425.9231

# Sheet LTW, row 40 (Leve Item ID 36248)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2990
$ws.Range("I40").Value = 2990
$ws.Range("K40").Value = 2990
$ws.Range("M40").Value = -2854

# Sheet LTW, row 93 (Leve Item ID 19993)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1666.3334
$ws.Range("J93").Value = 1666.3334
$ws.Range("L93").Value = 1666.3334
$ws.Range("N93").Value = -4162.3334

# Sheet LTW, row 100 (Leve Item ID 19995)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 4000
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").Value = ""

# Sheet LTW, row 122 (Leve Item ID 36247)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4095
$ws.Range("I122").Value = 2126.6667
$ws.Range("J122").Value = 10000
$ws.Range("K122").Value = 6380.000100000001
$ws.Range("L122").Value = 30000
$ws.Range("M122").Value = -3930.000100000001
$ws.Range("N122").Value = -34900

# Sheet LTW, row 132 (Leve Item ID 44058)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 7761.5
$ws.Range("I132").Value = 7441.7144
$ws.Range("K132").Value = 22325.1432
$ws.Range("M132").Value = -19795.1432

# Sheet WVR, row 109 (Leve Item ID 27161)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 18000
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").Value = ""

# Sheet WVR, row 113 (Leve Item ID 27752)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1782.5
$ws.Range("I113").Value = 2166.6667
$ws.Range("J113").Value = 1398.3334
$ws.Range("K113").Value = 6500.000100000001
$ws.Range("L113").Value = 4195.0002
$ws.Range("M113").Value = -4330.000100000001
$ws.Range("N113").Value = -8535.0002

# Sheet WVR, row 126 (Leve Item ID 36210)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2677
$ws.Range("I126").Value = 1281.6666
$ws.Range("K126").Value = 3844.9998
$ws.Range("M126").Value = -1374.9998
